# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" right before "总计" (mirrors the existing
#   per-quarter sheets' layout/format).
# - Insert a new top data-row in "总计" summarising the new quarter, and
#   shift the existing rows (and their running index in column A) down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet, positioned just before "总计".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q2")

$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Copy the header / row-2 cell formatting (borders, bold, alignment) from
# an existing quarter sheet so the new sheet matches the house style.
$template.Range("B1:H2").Copy()
$newSheet.Range("B1:H2").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Headers
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "580008"
$newSheet.Range("B2").Style = "Normal"

$newSheet.Range("C2").Value = "东吴新产业精选股票A"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "2.32"
$newSheet.Range("D2").Style = "Normal"

$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "89.77"
$newSheet.Range("E2").Style = "Normal"

$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "4.07"
$newSheet.Range("F2").Style = "Normal"

$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0944"
$newSheet.Range("G2").Style = "Normal"

$newSheet.Range("H2").Value = 7

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: add a 2022-Q1 row at the top of
#    the data and push the rest down, renumbering the index column.
#    (Re-fetch the sheet by name: worksheet references captured before
#    a Worksheets.Add() track sheet *position*, not identity, and
#    Add() shifts everything from the insertion point onward.)
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A2").Style = "Normal"
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("B2").Style = "Normal"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("C2").Style = "Normal"
$totalSheet.Range("D2").Value = 0.09
$totalSheet.Range("D2").Style = "Normal"

# Restore the row-2 "index" cell style (A column uses the bold/bordered
# style shared by the rest of the index column).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("A2").Value = 0

# Renumber the running index for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# ---------------------------------------------------------------------
# 3. Leave the original first sheet active, as in the source workbook.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
